# Update the marksheet's correct/total mark figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row (row 11): number of right answers used for marking -> 5
$ws.Range("B11").Value = 5

# "Total" row (row 12): total marks scored -> 110
$ws.Range("B12").Value = 110

# Score fraction text "corr/total" -> "110/140"
$ws.Range("E12").Value = "110/140"
